# Actualización automática 2025-11-17 16:30:09
$wb = $excel.ActiveWorkbook

$wsVentasPorGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsVentaMensual    = $wb.Worksheets.Item("VENTA MENSUAL")
$wsCumplimiento    = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# --- VENTAS POR GRUPO ---
$wsVentasPorGrupo.Range("M16").Value = 416.87
$wsVentasPorGrupo.Range("M21").Value = 1827.02

# --- VENTA MENSUAL ---
$wsVentaMensual.Range("F16").Value = 416.87
$wsVentaMensual.Range("F21").Value = 1850.42
$wsVentaMensual.Range("F38").Value = 2728.25

# --- CUMPLIMIENTO MENSUAL ---
$wsCumplimiento.Range("D12").Value = 3398.27
$wsCumplimiento.Range("E12").Value = 27815.73
$wsCumplimiento.Range("F12").Value = 0.1088700583071699

$wsCumplimiento.Range("D14").Value = 3487.91
$wsCumplimiento.Range("E14").Value = 36791.65164865473
$wsCumplimiento.Range("F14").Value = 0.08659255109139179
